# Add the "Growth from Previous Year" column (O) to the monthly-sales sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell O1 - same text style as the neighboring "Annual Sales" header (N1)
$ws.Range("O1").Value = "Growth from Previous Year"
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)  # xlPasteFormats

# O3 is a plain formula; O4:O11 share the same relative formula pattern
$ws.Range("O3").Formula = "=N3-N2"
$ws.Range("O4:O11").Formula = "=N4-N3"

# Give column O a sensible width so the header text is fully visible
$ws.Columns.Item(15).ColumnWidth = 22.7

# Match the author's final selection/active cell
$ws.Range("O17").Select()
